$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 50,7

$arr[0,0] = 1
$arr[0,1] = 'BTC'
$arr[0,2] = 'Bitcoin'
$arr[0,3] = 65481
$arr[0,4] = 1283567562392
$arr[0,5] = 24328168575
$arr[0,6] = 1.19805

$arr[1,0] = 2
$arr[1,1] = 'ETH'
$arr[1,2] = 'Ethereum'
$arr[1,3] = 3395.01
$arr[1,4] = 406561616932
$arr[1,5] = 13058337011
$arr[1,6] = 0.7223000000000001

$arr[2,0] = 3
$arr[2,1] = 'USDT'
$arr[2,2] = 'Tether'
$arr[2,3] = 1
$arr[2,4] = 103968330210
$arr[2,5] = 43973686919
$arr[2,6] = -0.10875

$arr[3,0] = 4
$arr[3,1] = 'BNB'
$arr[3,2] = 'BNB'
$arr[3,3] = 558.98
$arr[3,4] = 85768973301
$arr[3,5] = 1320270956
$arr[3,6] = 0.65806

$arr[4,0] = 5
$arr[4,1] = 'SOL'
$arr[4,2] = 'Solana'
$arr[4,3] = 175.55
$arr[4,4] = 77732934639
$arr[4,5] = 2723617530
$arr[4,6] = -0.43472

$arr[5,0] = 6
$arr[5,1] = 'XRP'
$arr[5,2] = 'XRP'
$arr[5,3] = 0.630345
$arr[5,4] = 34543249122
$arr[5,5] = 1206989490
$arr[5,6] = 1.70655

$arr[6,0] = 7
$arr[6,1] = 'STETH'
$arr[6,2] = 'Lido Staked Ether'
$arr[6,3] = 3385.43
$arr[6,4] = 33009276406
$arr[6,5] = 81706432
$arr[6,6] = 0.59617

$arr[7,0] = 8
$arr[7,1] = 'USDC'
$arr[7,2] = 'USDC'
$arr[7,3] = 1.002
$arr[7,4] = 32020837400
$arr[7,5] = 4767864913
$arr[7,6] = 0.0804

$arr[8,0] = 9
$arr[8,1] = 'DOGE'
$arr[8,2] = 'Dogecoin'
$arr[8,3] = 0.173545
$arr[8,4] = 24826396055
$arr[8,5] = 3077442030
$arr[8,6] = 5.36736

$arr[9,0] = 10
$arr[9,1] = 'ADA'
$arr[9,2] = 'Cardano'
$arr[9,3] = 0.637866
$arr[9,4] = 22431295701
$arr[9,5] = 377102053
$arr[9,6] = 1.1198

$arr[10,0] = 11
$arr[10,1] = 'AVAX'
$arr[10,2] = 'Avalanche'
$arr[10,3] = 53.8
$arr[10,4] = 20236798790
$arr[10,5] = 560032838
$arr[10,6] = -2.11694

$arr[11,0] = 12
$arr[11,1] = 'TON'
$arr[11,2] = 'Toncoin'
$arr[11,3] = 4.88
$arr[11,4] = 16895169186
$arr[11,5] = 387009199
$arr[11,6] = 1.43286

$arr[12,0] = 13
$arr[12,1] = 'SHIB'
$arr[12,2] = 'Shiba Inu'
$arr[12,3] = 0.00002804
$arr[12,4] = 16481980673
$arr[12,5] = 1005766188
$arr[12,6] = 1.61418

$arr[13,0] = 14
$arr[13,1] = 'DOT'
$arr[13,2] = 'Polkadot'
$arr[13,3] = 9.210000000000001
$arr[13,4] = 12362496107
$arr[13,5] = 175461629
$arr[13,6] = 1.18725

$arr[14,0] = 15
$arr[14,1] = 'LINK'
$arr[14,2] = 'Chainlink'
$arr[14,3] = 18.35
$arr[14,4] = 10753840142
$arr[14,5] = 280987063
$arr[14,6] = -0.00202

$arr[15,0] = 16
$arr[15,1] = 'TRX'
$arr[15,2] = 'TRON'
$arr[15,3] = 0.119023
$arr[15,4] = 10422556868
$arr[15,5] = 389083919
$arr[15,6] = 0.73532

$arr[16,0] = 17
$arr[16,1] = 'WBTC'
$arr[16,2] = 'Wrapped Bitcoin'
$arr[16,3] = 65408
$arr[16,4] = 10144753249
$arr[16,5] = 193352995
$arr[16,6] = 1.42038

$arr[17,0] = 18
$arr[17,1] = 'MATIC'
$arr[17,2] = 'Polygon'
$arr[17,3] = 1.002
$arr[17,4] = 9284027628
$arr[17,5] = 290446146
$arr[17,6] = 1.6188

$arr[18,0] = 19
$arr[18,1] = 'BCH'
$arr[18,2] = 'Bitcoin Cash'
$arr[18,3] = 463.21
$arr[18,4] = 9109806871
$arr[18,5] = 668722181
$arr[18,6] = 0.0153

$arr[19,0] = 20
$arr[19,1] = 'UNI'
$arr[19,2] = 'Uniswap'
$arr[19,3] = 11.87
$arr[19,4] = 8936120460
$arr[19,5] = 120044407
$arr[19,6] = 0.00101

$arr[20,0] = 21
$arr[20,1] = 'NEAR'
$arr[20,2] = 'NEAR Protocol'
$arr[20,3] = 6.56
$arr[20,4] = 6796340442
$arr[20,5] = 283000619
$arr[20,6] = -1.87724

$arr[21,0] = 22
$arr[21,1] = 'APT'
$arr[21,2] = 'Aptos'
$arr[21,3] = 16.66
$arr[21,4] = 6600529255
$arr[21,5] = 260320284
$arr[21,6] = 5.03464

$arr[22,0] = 23
$arr[22,1] = 'ICP'
$arr[22,2] = 'Internet Computer'
$arr[22,3] = 14.13
$arr[22,4] = 6532590555
$arr[22,5] = 246137727
$arr[22,6] = 5.68605

$arr[23,0] = 24
$arr[23,1] = 'LTC'
$arr[23,2] = 'Litecoin'
$arr[23,3] = 87.91
$arr[23,4] = 6525933928
$arr[23,5] = 539907098
$arr[23,6] = 1.82594

$arr[24,0] = 25
$arr[24,1] = 'LEO'
$arr[24,2] = 'LEO Token'
$arr[24,3] = 5.96
$arr[24,4] = 5584840656
$arr[24,5] = 1905622
$arr[24,6] = -1.84799

$arr[25,0] = 26
$arr[25,1] = 'STX'
$arr[25,2] = 'Stacks'
$arr[25,3] = 3.64
$arr[25,4] = 5269367761
$arr[25,5] = 150831143
$arr[25,6] = 4.59577

$arr[26,0] = 27
$arr[26,1] = 'DAI'
$arr[26,2] = 'Dai'
$arr[26,3] = 1
$arr[26,4] = 4898246257
$arr[26,5] = 632392419
$arr[26,6] = 0.18825

$arr[27,0] = 28
$arr[27,1] = 'FIL'
$arr[27,2] = 'Filecoin'
$arr[27,3] = 8.75
$arr[27,4] = 4601269851
$arr[27,5] = 206302815
$arr[27,6] = -0.25763

$arr[28,0] = 29
$arr[28,1] = 'ETC'
$arr[28,2] = 'Ethereum Classic'
$arr[28,3] = 31.19
$arr[28,4] = 4544586158
$arr[28,5] = 206777211
$arr[28,6] = 3.34493

$arr[29,0] = 30
$arr[29,1] = 'ATOM'
$arr[29,2] = 'Cosmos Hub'
$arr[29,3] = 11.5
$arr[29,4] = 4482624843
$arr[29,5] = 125783806
$arr[29,6] = -0.06662999999999999

$arr[30,0] = 31
$arr[30,1] = 'ARB'
$arr[30,2] = 'Arbitrum'
$arr[30,3] = 1.66
$arr[30,4] = 4391227721
$arr[30,5] = 332932527
$arr[30,6] = 2.76237

$arr[31,0] = 32
$arr[31,1] = 'IMX'
$arr[31,2] = 'Immutable'
$arr[31,3] = 2.9
$arr[31,4] = 4119167987
$arr[31,5] = 60584400
$arr[31,6] = 2.18426

$arr[32,0] = 33
$arr[32,1] = 'RNDR'
$arr[32,2] = 'Render'
$arr[32,3] = 10.72
$arr[32,4] = 4085042563
$arr[32,5] = 144292582
$arr[32,6] = -1.78188

$arr[33,0] = 34
$arr[33,1] = 'XLM'
$arr[33,2] = 'Stellar'
$arr[33,3] = 0.133973
$arr[33,4] = 3839965345
$arr[33,5] = 78114191
$arr[33,6] = 2.32075

$arr[34,0] = 35
$arr[34,1] = 'OKB'
$arr[34,2] = 'OKB'
$arr[34,3] = 63.53
$arr[34,4] = 3802093724
$arr[34,5] = 16424845
$arr[34,6] = 7.32737

$arr[35,0] = 36
$arr[35,1] = 'TAO'
$arr[35,2] = 'Bittensor'
$arr[35,3] = 581.2
$arr[35,4] = 3753773682
$arr[35,5] = 18636199
$arr[35,6] = -0.0633

$arr[36,0] = 37
$arr[36,1] = 'OP'
$arr[36,2] = 'Optimism'
$arr[36,3] = 3.69
$arr[36,4] = 3707640826
$arr[36,5] = 268843205
$arr[36,6] = 7.67278

$arr[37,0] = 38
$arr[37,1] = 'CRO'
$arr[37,2] = 'Cronos'
$arr[37,3] = 0.139415
$arr[37,4] = 3697585402
$arr[37,5] = 21707973
$arr[37,6] = 0.5129

$arr[38,0] = 39
$arr[38,1] = 'HBAR'
$arr[38,2] = 'Hedera'
$arr[38,3] = 0.108016
$arr[38,4] = 3631863698
$arr[38,5] = 39545337
$arr[38,6] = -0.3276

$arr[39,0] = 40
$arr[39,1] = 'GRT'
$arr[39,2] = 'The Graph'
$arr[39,3] = 0.373521
$arr[39,4] = 3521474857
$arr[39,5] = 114109920
$arr[39,6] = 0.16339

$arr[40,0] = 41
$arr[40,1] = 'KAS'
$arr[40,2] = 'Kaspa'
$arr[40,3] = 0.142686
$arr[40,4] = 3285669211
$arr[40,5] = 39935998
$arr[40,6] = 1.56118

$arr[41,0] = 42
$arr[41,1] = 'INJ'
$arr[41,2] = 'Injective'
$arr[41,3] = 35.85
$arr[41,4] = 3160342284
$arr[41,5] = 81952955
$arr[41,6] = -0.07027

$arr[42,0] = 43
$arr[42,1] = 'PEPE'
$arr[42,2] = 'Pepe'
$arr[42,3] = 0.00000745
$arr[42,4] = 3125220704
$arr[42,5] = 519655834
$arr[42,6] = -1.89398

$arr[43,0] = 44
$arr[43,1] = 'VET'
$arr[43,2] = 'VeChain'
$arr[43,3] = 0.04194042
$arr[43,4] = 3038933294
$arr[43,5] = 56060210
$arr[43,6] = 1.66337

$arr[44,0] = 45
$arr[44,1] = 'FTM'
$arr[44,2] = 'Fantom'
$arr[44,3] = 1.064
$arr[44,4] = 2948186270
$arr[44,5] = 315967790
$arr[44,6] = -3.46246

$arr[45,0] = 46
$arr[45,1] = 'MKR'
$arr[45,2] = 'Maker'
$arr[45,3] = 3099.58
$arr[45,4] = 2858305401
$arr[45,5] = 90077218
$arr[45,6] = -0.38831

$arr[46,0] = 47
$arr[46,1] = 'RUNE'
$arr[46,2] = 'THORChain'
$arr[46,3] = 8.369999999999999
$arr[46,4] = 2803228850
$arr[46,5] = 359728074
$arr[46,6] = -0.0606

$arr[47,0] = 48
$arr[47,1] = 'THETA'
$arr[47,2] = 'Theta Network'
$arr[47,3] = 2.78
$arr[47,4] = 2775639607
$arr[47,5] = 43253969
$arr[47,6] = -0.76181

$arr[48,0] = 49
$arr[48,1] = 'LDO'
$arr[48,2] = 'Lido DAO'
$arr[48,3] = 3.05
$arr[48,4] = 2680421105
$arr[48,5] = 112328792
$arr[48,6] = 7.22221

$arr[49,0] = 50
$arr[49,1] = 'MNT'
$arr[49,2] = 'Mantle'
$arr[49,3] = 0.825675
$arr[49,4] = 2668958705
$arr[49,5] = 133150251
$arr[49,6] = 0.38305

$ws.Range("A2:G51").Value = $arr
